$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prices")

$rowRange = $ws.Range("A76:J76")
$rowRange.NumberFormat = "@"

$ws.Range("A76").Value = "2025-05-16"
$ws.Range("B76").Value = "37.5"
$ws.Range("C76").Value = "37"
$ws.Range("D76").Value = "0.95"
$ws.Range("E76").Value = "0.258"
$ws.Range("F76").Value = "0.09"
$ws.Range("G76").Value = "5,259"
$ws.Range("H76").Value = "7,873"
$ws.Range("I76").Value = "7,923"
$ws.Range("J76").Value = "7.2286"

$rowRange.Style = "Normal"
